# This script applies a cyclic re-shuffling of the observation records held
# in rows 2-7 of the active worksheet. The header (row 1) and rows 8-14 are
# left untouched; only the "variable" columns (A, B, E, F, G, H, P, Q, R, AI,
# AO) that differ between the six records are rewritten, using the values
# that the record occupying each row acquires after the reshuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the "variable" columns of rows 2-7 after the edit.
# (Derived from rotating the original row 2-7 records: 2<-4, 3<-7, 4<-2,
# 5<-3, 6<-5, 7<-6.)

# Row 2 (was row 4's data)
$ws.Range("A2").Value = 865900
$ws.Range("B2").Value = 78569
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("P2").Value = "SV Monte Carlo, Jmt"
$ws.Range("Q2").Value = 527721.6135952856
$ws.Range("R2").Value = 6996047.266259992
$ws.Range("AI2").ClearContents()
$ws.Range("AO2").Value = "1 substratenheter # sälgbark"

# Row 3 (was row 7's data)
$ws.Range("A3").Value = 1955206
$ws.Range("B3").Value = 78568
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("P3").Value = "SV Monte Carlo, Jmt"
$ws.Range("Q3").Value = 527646.3222403944
$ws.Range("R3").Value = 6996046.623740066
$ws.Range("AI3").ClearContents()
$ws.Range("AO3").Value = "1 substratenheter # björkbark"

# Row 4 (was row 2's data)
$ws.Range("A4").Value = 1955211
$ws.Range("B4").Value = 78568
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("P4").Value = "Söder om Bodtjärnen, Jmt"
$ws.Range("Q4").Value = 526984.7824697205
$ws.Range("R4").Value = 6996343.169846137
$ws.Range("AI4").Value = "grannaturskog"
$ws.Range("AO4").Value = "1 substratenheter # gammal björk"

# Row 5 (was row 3's data)
$ws.Range("A5").Value = 1955208
$ws.Range("B5").Value = 78568
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("P5").Value = "SV Monte Carlo, Jmt"
$ws.Range("Q5").Value = 527721.6135952856
$ws.Range("R5").Value = 6996047.266259992
$ws.Range("AO5").Value = "1 substratenheter"

# Row 6 (was row 5's data)
$ws.Range("A6").Value = 1955205
$ws.Range("B6").Value = 78568
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("P6").Value = "S Monte Carlo, Jmt"
$ws.Range("Q6").Value = 527801.0263964261
$ws.Range("R6").Value = 6995990.418098574
$ws.Range("AO6").Value = "1 substratenheter # aspbark"

# Row 7 (was row 6's data)
$ws.Range("A7").Value = 1955210
$ws.Range("B7").Value = 78568
$ws.Range("E7").Value = 6458
$ws.Range("F7").Value = "Lunglav"
$ws.Range("G7").Value = "Lobaria pulmonaria"
$ws.Range("H7").Value = "(L.) Hoffm."
$ws.Range("P7").Value = "Söder om Bodtjärnen, Jmt"
$ws.Range("Q7").Value = 527394.2661032595
$ws.Range("R7").Value = 6996083.440633372
$ws.Range("AI7").Value = "grannaturskog"
$ws.Range("AO7").Value = "1 substratenheter # björk med brandljud"
